# Updated symbol list on Fri Jan 13 04:13:00 UTC 2023 with GitHub Actions
# Refresh of coin price/volume/hour data in the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.65%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.89%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "4"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.141"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.48%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06667"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.78%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.351"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.50%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.404"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.43%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "4"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.80%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9161"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.04%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1580"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.68%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "4"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06611"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.11%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07646"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.14%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02974"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.22%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "4"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08999"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "4"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001594"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.41%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "4"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04473"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.07%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "4"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006454"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.29%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "4"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006262"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.51%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "4"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.458"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.86%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.76%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3214"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.81%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "4"
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.93%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "4"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.074"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.35%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "4"
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.03%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "4"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001189"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.09%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "4"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004138"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.65%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "4"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001247"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.78%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "4"
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.16%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "4"
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "4"
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "4"
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "4"
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "4"
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "4"
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "4"
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "4"
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "4"
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "4"
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "4"
# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "4"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04224"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.38%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "4"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006734"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.80%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "4"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1242"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-12.33%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "4"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001976"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.02%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "4"
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.74%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "4"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005599"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.90%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "4"
# Row 46
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.01305"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-29.42%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "4"
# Row 47
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.968"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "25.94%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "4"
# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "4"
# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "4"
# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "4"
# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "4"
